# Advance all timestamps and lookup date labels in Sheet1 by one day
# (data refresh / retrain window roll-forward)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newA = @(46047.99444444444,46048.00486111111,46048.01527777778,46048.02569444444,46048.03611111111,46048.04652777778,46048.05694444444,46048.06736111111,46048.07777777778,46048.08819444444,46048.09861111111,46048.10902777778,46048.11944444444,46048.12986111111,46048.14027777778,46048.15069444444,46048.16111111111,46048.17152777778,46048.18194444444,46048.19236111111,46048.20277777778,46048.21319444444,46048.22361111111,46048.23402777778,46048.24444444444,46048.25486111111,46048.26527777778,46048.27569444444,46048.28611111111,46048.29652777778,46048.30694444444,46048.31736111111,46048.32777777778,46048.33819444444,46048.34861111111,46048.35902777778,46048.36944444444,46048.37986111111,46048.39027777778,46048.40069444444,46048.41111111111,46048.42152777778,46048.43194444444,46048.44236111111,46048.45277777778,46048.46319444444,46048.47361111111,46048.48402777778,46048.49444444444,46048.50486111111,46048.51527777778,46048.52569444444,46048.53611111111,46048.54652777778,46048.55694444444,46048.56736111111,46048.57777777778,46048.58819444444,46048.59861111111,46048.60902777778,46048.61944444444,46048.62986111111,46048.64027777778,46048.65069444444,46048.66111111111,46048.67152777778,46048.68194444444,46048.69236111111,46048.70277777778,46048.71319444444,46048.72361111111,46048.73402777778,46048.74444444444,46048.75486111111,46048.76527777778,46048.77569444444,46048.78611111111,46048.79652777778,46048.80694444444,46048.81736111111,46048.82777777778,46048.83819444444,46048.84861111111,46048.85902777778,46048.86944444444,46048.87986111111,46048.89027777778,46048.90069444444,46048.91111111111,46048.92152777778,46048.93194444444,46048.94236111111,46048.95277777778,46048.96319444444,46048.97361111111,46048.98402777778,46048.99444444444,46048.99444444444,46049.00486111111,46049.01527777778,46049.02569444444,46049.03611111111,46049.04652777778,46049.05694444444,46049.06736111111,46049.07777777778,46049.08819444444,46049.09861111111,46049.10902777778,46049.11944444444,46049.12986111111,46049.14027777778,46049.15069444444,46049.16111111111,46049.17152777778,46049.18194444444,46049.19236111111,46049.20277777778,46049.21319444444,46049.22361111111,46049.23402777778,46049.24444444444,46049.25486111111,46049.26527777778,46049.27569444444,46049.28611111111,46049.29652777778,46049.30694444444,46049.31736111111,46049.32777777778,46049.33819444444,46049.34861111111,46049.35902777778,46049.36944444444,46049.37986111111,46049.39027777778,46049.40069444444,46049.41111111111,46049.42152777778,46049.43194444444,46049.44236111111,46049.45277777778,46049.46319444444,46049.47361111111,46049.48402777778,46049.49444444444,46049.50486111111,46049.51527777778,46049.52569444444,46049.53611111111,46049.54652777778,46049.55694444444,46049.56736111111,46049.57777777778,46049.58819444444,46049.59861111111,46049.60902777778,46049.61944444444,46049.62986111111,46049.64027777778,46049.65069444444,46049.66111111111,46049.67152777778,46049.68194444444,46049.69236111111,46049.70277777778,46049.71319444444,46049.72361111111,46049.73402777778,46049.74444444444,46049.75486111111,46049.76527777778,46049.77569444444,46049.78611111111,46049.79652777778,46049.80694444444,46049.81736111111,46049.82777777778,46049.83819444444,46049.84861111111,46049.85902777778,46049.86944444444,46049.87986111111,46049.89027777778,46049.90069444444,46049.91111111111,46049.92152777778,46049.93194444444,46049.94236111111,46049.95277777778,46049.96319444444,46049.97361111111,46049.98402777778,46049.99444444444)
$newE = @("25.01.20261","26.01.20262","26.01.20263","26.01.20264","26.01.20265","26.01.20266","26.01.20267","26.01.20268","26.01.20269","26.01.202610","26.01.202611","26.01.202612","26.01.202613","26.01.202614","26.01.202615","26.01.202616","26.01.202617","26.01.202618","26.01.202619","26.01.202620","26.01.202621","26.01.202622","26.01.202623","26.01.202624","26.01.202625","26.01.202626","26.01.202627","26.01.202628","26.01.202629","26.01.202630","26.01.202631","26.01.202632","26.01.202633","26.01.202634","26.01.202635","26.01.202636","26.01.202637","26.01.202638","26.01.202639","26.01.202640","26.01.202641","26.01.202642","26.01.202643","26.01.202644","26.01.202645","26.01.202646","26.01.202647","26.01.202648","26.01.202649","26.01.202650","26.01.202651","26.01.202652","26.01.202653","26.01.202654","26.01.202655","26.01.202656","26.01.202657","26.01.202658","26.01.202659","26.01.202660","26.01.202661","26.01.202662","26.01.202663","26.01.202664","26.01.202665","26.01.202666","26.01.202667","26.01.202668","26.01.202669","26.01.202670","26.01.202671","26.01.202672","26.01.202673","26.01.202674","26.01.202675","26.01.202676","26.01.202677","26.01.202678","26.01.202679","26.01.202680","26.01.202681","26.01.202682","26.01.202683","26.01.202684","26.01.202685","26.01.202686","26.01.202687","26.01.202688","26.01.202689","26.01.202690","26.01.202691","26.01.202692","26.01.202693","26.01.202694","26.01.202695","26.01.202696","26.01.20261","26.01.20262","27.01.20263","27.01.20264","27.01.20265","27.01.20266","27.01.20267","27.01.20268","27.01.20269","27.01.202610","27.01.202611","27.01.202612","27.01.202613","27.01.202614","27.01.202615","27.01.202616","27.01.202617","27.01.202618","27.01.202619","27.01.202620","27.01.202621","27.01.202622","27.01.202623","27.01.202624","27.01.202625","27.01.202626","27.01.202627","27.01.202628","27.01.202629","27.01.202630","27.01.202631","27.01.202632","27.01.202633","27.01.202634","27.01.202635","27.01.202636","27.01.202637","27.01.202638","27.01.202639","27.01.202640","27.01.202641","27.01.202642","27.01.202643","27.01.202644","27.01.202645","27.01.202646","27.01.202647","27.01.202648","27.01.202649","27.01.202650","27.01.202651","27.01.202652","27.01.202653","27.01.202654","27.01.202655","27.01.202656","27.01.202657","27.01.202658","27.01.202659","27.01.202660","27.01.202661","27.01.202662","27.01.202663","27.01.202664","27.01.202665","27.01.202666","27.01.202667","27.01.202668","27.01.202669","27.01.202670","27.01.202671","27.01.202672","27.01.202673","27.01.202674","27.01.202675","27.01.202676","27.01.202677","27.01.202678","27.01.202679","27.01.202680","27.01.202681","27.01.202682","27.01.202683","27.01.202684","27.01.202685","27.01.202686","27.01.202687","27.01.202688","27.01.202689","27.01.202690","27.01.202691","27.01.202692","27.01.202693","27.01.202694","27.01.202695","27.01.202696","27.01.20261","27.01.20262")

$startRow = 2
for ($i = 0; $i -lt $newA.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $newA[$i]
    $ws.Cells.Item($r, 5).Value2 = $newE[$i]
}

"Done: updated " + $newA.Length + " rows"
